# Append the 09/14/2025 profit-allocation row (row 13) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format as Text first so the "MM/DD/YYYY"-shaped string is kept as a literal
# string instead of being auto-converted into a date serial number, then
# reset the style back to Normal so the cell ends up with no explicit style
# (matching the other date cells in column A, which are plain strings).
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "09/14/2025"
$ws.Range("A13").Style = "Normal"

$ws.Range("B13").Value = 0.1222829562546641
$ws.Range("C13").Value = 0.8777170437453359
